$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "PASSED"
$ws.Range("E2").Value = 8.212997483002255
$ws.Range("F2").Value = "2022-08-19T17:20:39"
$ws.Range("G2").ClearContents()
